$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row: "Ahmedabad" location assertion (rows 1-2 already hold the
# right data - Assertions/Locations/Expected header and the Delhi example -
# we only need to append the Ahmedabad example in row 3).
$ws.Range("B3").Value = "Ahmedabad"
$ws.Range("C3").Value = "Ahmedabad, Gujarat, India"

# Column B is brand-new (holds "Locations"/"Delhi"/"Ahmedabad"); give it a
# sensible best-fit-ish width. Column C got wider to fit the longer
# "Ahmedabad, Gujarat, India" expected text.
$ws.Columns.Item(2).ColumnWidth = 10.75
$ws.Columns.Item(3).ColumnWidth = 24.1

# Selection moved to C3 (the newly entered cell) before the file was saved.
$ws.Range("C3").Select()
